$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 62288.83
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 64106.15
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 192318.45
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -192654.45
# Row 51
$ws.Range("H51").Value = 11113017
$ws.Range("I51").Value = 1923.0769
$ws.Range("J51").Value = 40001860
$ws.Range("K51").Value = 1923.0769
$ws.Range("L51").Value = 40001860
$ws.Range("M51").Value = -1439.0769
$ws.Range("N51").Value = -40002828
# Row 86
$ws.Range("H86").Value = 3114.25
$ws.Range("I86").Value = 2226.5
$ws.Range("K86").Value = 2226.5
$ws.Range("M86").Value = -1103.5
# Row 89
$ws.Range("H89").Value = 3114.25
$ws.Range("I89").Value = 2226.5
$ws.Range("K89").Value = 11132.5
$ws.Range("M89").Value = -5516.5
# Row 100
$ws.Range("H100").Value = 2764.7727
$ws.Range("I100").Value = 2180.3125
$ws.Range("J100").Value = 4323.3335
$ws.Range("K100").Value = 2180.3125
$ws.Range("L100").Value = 4323.3335
$ws.Range("M100").Value = -1639.3125
$ws.Range("N100").Value = -5405.3335
# Row 111
$ws.Range("H111").Value = 765.5
$ws.Range("I111").Value = 660.2857
$ws.Range("K111").Value = 1980.8571
$ws.Range("M111").Value = 1086.1429

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 8669
$ws.Range("I61").Value = 1000
$ws.Range("K61").Value = 1000
$ws.Range("M61").Value = -788
# Row 63
$ws.Range("H63").Value = 2577
$ws.Range("I63").Value = 1884.9286
$ws.Range("K63").Value = 1884.9286
$ws.Range("M63").Value = -1198.9286
# Row 66
$ws.Range("H66").Value = 2577
$ws.Range("I66").Value = 1884.9286
$ws.Range("K66").Value = 9424.643
$ws.Range("M66").Value = -5992.643
# Row 74
$ws.Range("I74").Value = 1770.3334
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1770.3334
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -896.3334
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("I77").Value = 1770.3334
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 8851.666999999999
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -4483.666999999999
$ws.Range("N77").ClearContents()
# Row 110
$ws.Range("H110").Value = 2045.1428
$ws.Range("I110").Value = 481.9
$ws.Range("J110").Value = 5953.25
$ws.Range("K110").Value = 481.9
$ws.Range("L110").Value = 5953.25
$ws.Range("M110").Value = 1563.1
$ws.Range("N110").Value = -10043.25
# Row 136
$ws.Range("H136").Value = 8669
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 17
$ws.Range("H17").Value = 2009
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2009
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2009
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -2353

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 32717.143
$ws.Range("I2").Value = 3000
$ws.Range("K2").Value = 3000
$ws.Range("M2").Value = -2887
# Row 10
$ws.Range("H10").Value = 7404.92
$ws.Range("I10").Value = 760.4666999999999
$ws.Range("K10").Value = 760.4666999999999
$ws.Range("M10").Value = -621.4666999999999
# Row 122
$ws.Range("H122").Value = 1636.0541
$ws.Range("I122").Value = 1402.3462
$ws.Range("J122").Value = 2188.4546
$ws.Range("K122").Value = 4207.0386
$ws.Range("L122").Value = 6565.3638
$ws.Range("M122").Value = -1757.0386
$ws.Range("N122").Value = -11465.3638
# Row 134
$ws.Range("H134").Value = 3068.3076
$ws.Range("I134").Value = 1188.8
$ws.Range("K134").Value = 3566.4
$ws.Range("M134").Value = -1031.4

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 18
$ws.Range("H18").Value = 722
$ws.Range("I18").Value = 236.66667
$ws.Range("J18").Value = 1450
$ws.Range("K18").Value = 710.00001
$ws.Range("L18").Value = 4350
$ws.Range("M18").Value = -541.00001
$ws.Range("N18").Value = -4688
# Row 131
$ws.Range("H131").Value = 1453.8108
$ws.Range("J131").Value = 1184.4412
$ws.Range("L131").Value = 3553.3236
$ws.Range("N131").Value = -13633.3236

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = 25626.5
$ws.Range("I19").Value = 4500
$ws.Range("J19").Value = 32668.666
$ws.Range("K19").Value = 4500
$ws.Range("L19").Value = 32668.666
$ws.Range("M19").Value = -4212
$ws.Range("N19").Value = -33244.666
# Row 80
$ws.Range("H80").Value = 2985.7144
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 2900
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 2900
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -4896
# Row 83
$ws.Range("H83").Value = 2985.7144
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 2900
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 14500
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -24484
# Row 102
$ws.Range("H102").Value = 41353.848
$ws.Range("I102").Value = 2327.4
$ws.Range("J102").Value = 94571.73
$ws.Range("K102").Value = 2327.4
$ws.Range("L102").Value = 94571.73
$ws.Range("M102").Value = -705.4000000000001
$ws.Range("N102").Value = -97815.73
# Row 107
$ws.Range("H107").Value = 1363.1818
$ws.Range("I107").Value = 356.42856
$ws.Range("J107").Value = 3125
$ws.Range("K107").Value = 356.42856
$ws.Range("L107").Value = 3125
$ws.Range("M107").Value = 1563.57144
$ws.Range("N107").Value = -6965
# Row 122
$ws.Range("H122").Value = 4946.706
$ws.Range("I122").Value = 4169
$ws.Range("K122").Value = 12507
$ws.Range("M122").Value = -10057
# Row 126
$ws.Range("H126").Value = 3938.375
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
# Row 132
$ws.Range("H132").Value = 4645.7036
$ws.Range("I132").Value = 5388.2856
$ws.Range("K132").Value = 16164.8568
$ws.Range("M132").Value = -13634.8568

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 7709.778
$ws.Range("I40").Value = 10677.6
$ws.Range("K40").Value = 10677.6
$ws.Range("M40").Value = -10541.6
# Row 46
$ws.Range("H46").Value = 1408.909
$ws.Range("I46").Value = 909.6799999999999
$ws.Range("J46").Value = 2969
$ws.Range("K46").Value = 909.6799999999999
$ws.Range("L46").Value = 909.6799999999999
$ws.Range("M46").Value = -721.6799999999999
$ws.Range("N46").Value = -3345
# Row 122
$ws.Range("H122").Value = 2687.7896
$ws.Range("I122").Value = 2320.3076
$ws.Range("K122").Value = 6960.9228
$ws.Range("M122").Value = -4510.9228

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 2417999.8
$ws.Range("J3").Value = 44999.5
$ws.Range("L3").Value = 44999.5
$ws.Range("N3").Value = -45227.5
# Row 122
$ws.Range("H122").Value = 771900.0600000001
$ws.Range("I122").Value = 1430770.2
$ws.Range("J122").Value = 3218.1667
$ws.Range("K122").Value = 4292310.6
$ws.Range("L122").Value = 9654.500100000001
$ws.Range("M122").Value = -4289860.6
$ws.Range("N122").Value = -14554.5001
